# Fix duplicate num val: insert the missing G8/G9/G10/G11 "ports" rows,
# interleaved with the existing F8/F9/F10/F11 rows. This pushes the
# existing "enum_list" rows further down the sheet (their content is
# unaffected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the 4 new "ports" rows, top to bottom -------------------------
# Each insert pushes everything below it down by one row, so the target
# row numbers below already account for the rows inserted before them.

$ws.Rows("11:11").Insert()
$ws.Range("A11").Value = "ports"
$ws.Range("B11").Value = "G8"
$ws.Range("C11").Value = "drpa"
$ws.Range("D11").Value = "nan"

$ws.Rows("13:13").Insert()
$ws.Range("A13").Value = "ports"
$ws.Range("B13").Value = "G9"
$ws.Range("C13").Value = "drpv"
$ws.Range("D13").Value = "nan"

$ws.Rows("15:15").Insert()
$ws.Range("A15").Value = "ports"
$ws.Range("B15").Value = "G10"
$ws.Range("C15").Value = "dra"
$ws.Range("D15").Value = "nan"

$ws.Rows("17:17").Insert()
$ws.Range("A17").Value = "ports"
$ws.Range("B17").Value = "G11"
$ws.Range("C17").Value = "drpa"
$ws.Range("D17").Value = "nan"

# The pre-existing "enum_list" rows (previously rows 14-23) are pushed down
# by the four inserts above and land at rows 18-27 with their values
# untouched -- no further changes are required there.
